$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set header cell C1, copying format from B1 (bold, border, centered) then set the text
$ws.Range("B1").Copy()
$ws.Range("C1").PasteSpecial(-4122)
$ws.Range("C1").Value = "Federal Contract Compensation"

# Federal Contract Compensation values for rows 2-84 (County rows)
$values = @(
    147820.5,
    215157.64,
    0,
    1370792.5,
    61639375.49,
    6483870.44,
    3231566.22,
    0,
    0,
    86482253.62,
    380477,
    0,
    0,
    0,
    3600.99,
    111280345.73,
    79663,
    408668.69,
    0,
    8155410.51,
    96224475.73999999,
    12075683.78,
    0,
    349409.25,
    7607370.51,
    568311.13,
    4226963.22,
    228754.85,
    7641129.11,
    38600.28,
    42547594.5,
    65141.88,
    1417.4,
    0,
    12825763.36,
    650467.25,
    55100,
    514166.44,
    53758188.34,
    0,
    33374734.17,
    0,
    0,
    240655,
    0,
    1216175.75,
    0,
    78487432.45,
    0,
    1823625.13,
    89726.92999999999,
    5243408.53,
    4511160.52,
    0,
    7224392.34,
    0,
    0,
    0,
    17276567.55,
    309009.69,
    55776518.4,
    0,
    58185439.92,
    8165960.16,
    108452.7,
    8200,
    12264562.82,
    45064.24,
    347395.52,
    6975560.97,
    865034.62,
    0,
    11686,
    1424479.58,
    357106,
    192024.89,
    0,
    70737.34,
    12619245.46,
    0,
    6525493.61,
    208230619.67,
    5772160.6
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 3).Value = $values[$i]
}

